$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.260331749916077
$ws.Range("B1").Value = 2.527748346328735
$ws.Range("C1").Value = 3.661454200744629
$ws.Range("D1").Value = 2.887630701065063
$ws.Range("E1").Value = 1.070904493331909
